# Auto-generated edit script: updates crypto price/volume table
# per commit "Updated cryptos list on Sun Sep 22 07:45:38 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.922.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.583.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.045.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.795.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.587.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.708.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.37%  "
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "468.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.77%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0821"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "175.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "157.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.634"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0540"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0965"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
